$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("2024-03-21 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-22 Friday", 2)
$null = $d.Content.Find.Execute("583×3=1749", $true, $false, $false, $false, $false, $true, 1, $false, "204×3=612", 2)
$null = $d.Content.Find.Execute("599×4=2396", $true, $false, $false, $false, $false, $true, 1, $false, "509×5=2545", 2)
$null = $d.Content.Find.Execute("300×7=2100", $true, $false, $false, $false, $false, $true, 1, $false, "268×4=1072", 2)
$null = $d.Content.Find.Execute("594×9=5346", $true, $false, $false, $false, $false, $true, 1, $false, "537×4=2148", 2)
$null = $d.Content.Find.Execute("290×4=1160", $true, $false, $false, $false, $false, $true, 1, $false, "993×2=1986", 2)
$null = $d.Content.Find.Execute("945×3=2835", $true, $false, $false, $false, $false, $true, 1, $false, "492×8=3936", 2)
$null = $d.Content.Find.Execute("115×7=805", $true, $false, $false, $false, $false, $true, 1, $false, "767×5=3835", 2)
$null = $d.Content.Find.Execute("781×5=3905", $true, $false, $false, $false, $false, $true, 1, $false, "424×4=1696", 2)
$null = $d.Content.Find.Execute("116×9=1044", $true, $false, $false, $false, $false, $true, 1, $false, "797×6=4782", 2)
$null = $d.Content.Find.Execute("697×6=4182", $true, $false, $false, $false, $false, $true, 1, $false, "441×5=2205", 2)
$null = $d.Content.Find.Execute("733×4=2932", $true, $false, $false, $false, $false, $true, 1, $false, "607×2=1214", 2)
$null = $d.Content.Find.Execute("440×7=3080", $true, $false, $false, $false, $false, $true, 1, $false, "265×3=795", 2)
$null = $d.Content.Find.Execute("367×2=734", $true, $false, $false, $false, $false, $true, 1, $false, "720×4=2880", 2)
$null = $d.Content.Find.Execute("354×3=1062", $true, $false, $false, $false, $false, $true, 1, $false, "473×8=3784", 2)
$null = $d.Content.Find.Execute("213×5=1065", $true, $false, $false, $false, $false, $true, 1, $false, "299×4=1196", 2)
$null = $d.Content.Find.Execute("931×7=6517", $true, $false, $false, $false, $false, $true, 1, $false, "843×4=3372", 2)
$null = $d.Content.Find.Execute("813×4=3252", $true, $false, $false, $false, $false, $true, 1, $false, "148×3=444", 2)
$null = $d.Content.Find.Execute("239×8=1912", $true, $false, $false, $false, $false, $true, 1, $false, "194×7=1358", 2)
$null = $d.Content.Find.Execute("918×5=4590", $true, $false, $false, $false, $false, $true, 1, $false, "393×8=3144", 2)
$null = $d.Content.Find.Execute("708×4=2832", $true, $false, $false, $false, $false, $true, 1, $false, "492×3=1476", 2)
$null = $d.Content.Find.Execute("395×7=2765", $true, $false, $false, $false, $false, $true, 1, $false, "930×7=6510", 2)
$null = $d.Content.Find.Execute("115×6=690", $true, $false, $false, $false, $false, $true, 1, $false, "692×9=6228", 2)
$null = $d.Content.Find.Execute("410×8=3280", $true, $false, $false, $false, $false, $true, 1, $false, "348×6=2088", 2)
$null = $d.Content.Find.Execute("333×4=1332", $true, $false, $false, $false, $false, $true, 1, $false, "388×9=3492", 2)
$null = $d.Content.Find.Execute("937×7=6559", $true, $false, $false, $false, $false, $true, 1, $false, "668×4=2672", 2)
